$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet is protected; unprotect to allow writing to the locked data cells,
# then restore protection afterwards.
$ws.Unprotect()

# Update the "as of" date in the confidential disclaimer text (shared string).
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-23 for illustrative purposes only and are subject to change."

# Row 2 (NULG)
$ws.Range("D2").Value = 0.2470344001342337
$ws.Range("E2").Value = -0.003150708909504751

# Row 3 (NULV)
$ws.Range("D3").Value = 0.4942059467453053
$ws.Range("E3").Value = -0.01147174034695031

# Row 4 (NUMG)
$ws.Range("D4").Value = 0.09939511790078881
$ws.Range("E4").Value = -0.01086956521739124

# Row 5 (NUMV)
$ws.Range("D5").Value = 0.1007129946063734
$ws.Range("E5").Value = -0.02141203703703709

# Row 6 (NUSC)
$ws.Range("D6").Value = 0.05865154061329884
$ws.Range("E6").Value = -0.03265118418027146

# Row 7 (Total) - only Percent Change changes
$ws.Range("E7").Value = -0.01159963012639487

# Restore sheet protection (legacy hashed password "D382", matching the
# original file's protection record).
$ws.Protect("D382")
